# QR Game Template v500C edit script
# 1. Remove the stray "_GoBack" bookmark that sat after "will delete them "
#    in the body of the document.
# 2. In the title-page (first page) header, turn the 3-run hyperlink field
#    https://qrproblems.org/crud into a single plain (non-linked) run reading
#    https://qrproblems.org, and drop a fresh "_GoBack" bookmark right after it
#    (this is where Word now leaves the cursor / last-edit marker).
# 3. Turn on distinct odd/even headers & footers for the section, which causes
#    Word to materialize header1/2/3.xml + footer1/2/3.xml parts and rewire the
#    section's header/footer references (even/default/first).

$d = $word.ActiveDocument

# --- 1. drop the old _GoBack bookmark in the body -------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. fix the QR-problems URL in the first-page header -------------------
$sec = $d.Sections(1)
$firstHeader = $sec.Headers(2)   # wdHeaderFooterFirstPage

if ($firstHeader.Range.Hyperlinks.Count -gt 0) {
    $firstHeader.Range.Hyperlinks.Item(1).Delete()
}

$hdrRange = $firstHeader.Range
$hdrRange.Find.Execute("https://qrproblems.org/crud", $true, $false, $false, $false, $false, $true, 1, $false, "https://qrproblems.org", 2) | Out-Null
$hdrRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $hdrRange) | Out-Null

# --- 3. enable distinct even/first headers & footers ------------------------
$sec.PageSetup.OddAndEvenPagesHeaderFooter = -1
$sec.Headers(1).Range.Text = ""
$sec.Headers(3).Range.Text = ""
$sec.Footers(1).Range.Text = ""
$sec.Footers(2).Range.Text = ""
$sec.Footers(3).Range.Text = ""
